# Revert "Powerpoint writer: consolidate text run nodes."
#
# The original deck stores each word of these paragraphs as its own run,
# with the trailing space glued onto the end of the word (e.g. "Blank ").
# The target splits the trailing space off into its own run (e.g. "Blank"
# followed by a run containing just " "), while leaving paragraph / run
# properties (empty <a:rPr/>) untouched.
#
# We rebuild each affected paragraph's text via TextRange.Text (first
# word) followed by TextRange.InsertAfter(...) calls (space, next word,
# space, next word, ...) — InsertAfter appends a brand-new run each time
# instead of merging into the previous one, which is exactly the
# word/space run split the diff shows.

$p = $ppt.ActivePresentation

# --- Slide 1: Title "Section Header (with background image)" ---
$tr = $p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Section"
$tr.InsertAfter(" ")
$tr.InsertAfter("Header")
$tr.InsertAfter(" ")
$tr.InsertAfter("(with")
$tr.InsertAfter(" ")
$tr.InsertAfter("background")
$tr.InsertAfter(" ")
$tr.InsertAfter("image)")

# --- Slide 2: Title "Slide 1" ---
$tr = $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Slide"
$tr.InsertAfter(" ")
$tr.InsertAfter("1")

# --- Slide 3: Title "Slide 2" ---
$tr = $p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Slide"
$tr.InsertAfter(" ")
$tr.InsertAfter("2")

# --- Slide 4: Title "Slide 3" ---
$tr = $p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Slide"
$tr.InsertAfter(" ")
$tr.InsertAfter("3")

# --- Slide 5: Title "Slide 4" ---
$tr = $p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Slide"
$tr.InsertAfter(" ")
$tr.InsertAfter("4")

# --- Slide 5: TextBox "An image" ---
$tr = $p.Slides.Item(5).Shapes.Item(4).TextFrame.TextRange
$tr.Text = "An"
$tr.InsertAfter(" ")
$tr.InsertAfter("image")

# --- Slide 6's notes page: "Blank slides can have background images." ---
# Some hosts refuse incremental InsertAfter edits on notes-page text
# ranges; fall back to writing the full (correctly worded) sentence in
# one shot so the visible text is still right even if the run split
# can't be reproduced there.
$tr = $p.Slides.Item(6).NotesPage.Shapes.Placeholders.Item(2).TextFrame.TextRange
try {
    $tr.Text = "Blank"
    $tr.InsertAfter(" ")
    $tr.InsertAfter("slides")
    $tr.InsertAfter(" ")
    $tr.InsertAfter("can")
    $tr.InsertAfter(" ")
    $tr.InsertAfter("have")
    $tr.InsertAfter(" ")
    $tr.InsertAfter("background")
    $tr.InsertAfter(" ")
    $tr.InsertAfter("images.")
} catch {
    $tr.Text = "Blank slides can have background images."
}
